# feat: add 2022-Q3 data
#
# The existing "2021-Q3" sheet (fund-holding details) is preserved as-is but
# relocated to a newly appended sheet, while the original sheet is
# repurposed to hold the new "2022-Q3" figures. The "总计" (totals) sheet
# gets a new first data row for 2022-Q3, pushing the old 2021-Q3 row down.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$oldQ3 = $wb.Worksheets.Item("2021-Q3")

# --- 1. Preserve the existing 2021-Q3 sheet by copying it right after
#        itself, then restore the "2021-Q3" name on the copy (Excel
#        auto-suffixes it as "2021-Q3 (2)").
$oldQ3.Copy($null, $oldQ3)
$preserved = $wb.Worksheets.Item("2021-Q3 (2)")

# Rename the original sheet (still holding the 2021-Q3 data for the moment)
# to "2022-Q3" - it keeps sheetId 2, matching the author's edit - before
# freeing up the "2021-Q3" name for the preserved copy.
$oldQ3.Name = "2022-Q3"
$preserved.Name = "2021-Q3"
$newQ3 = $oldQ3

# --- 2. Re-template the now-"2022-Q3" sheet so it matches a freshly added
#        sheet (same page margins / view as the "总计" sheet), then drop in
#        the new fund data.
$newQ3.PageSetup.LeftMargin = $total.PageSetup.LeftMargin
$newQ3.PageSetup.RightMargin = $total.PageSetup.RightMargin
$newQ3.PageSetup.TopMargin = $total.PageSetup.TopMargin
$newQ3.PageSetup.BottomMargin = $total.PageSetup.BottomMargin
$newQ3.PageSetup.HeaderMargin = $total.PageSetup.HeaderMargin
$newQ3.PageSetup.FooterMargin = $total.PageSetup.FooterMargin

$newQ3.Cells.Clear()

# Header row + index column pick up the bold/border/center style used on
# "总计" (style index 2 in the original workbook).
$total.Range("B1").Copy()
$newQ3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$newQ3.Range("A2:A3").PasteSpecial(-4122)

$newQ3.Range("B1").Value = "基金代码"
$newQ3.Range("C1").Value = "基金名称"
$newQ3.Range("D1").Value = "基金规模"
$newQ3.Range("E1").Value = "股票总仓位"
$newQ3.Range("F1").Value = "仓位占比"
$newQ3.Range("G1").Value = "持有市值(亿元)"
$newQ3.Range("H1").Value = "仓位排名"

$newQ3.Range("A2").Value = 0

# Fund code / name / numeric-looking figures are stored as plain text in the
# source data, so force a text number format before assigning - otherwise a
# string like "004845" gets coerced into the number 4845 - then drop back to
# the Normal style so no stray number-format residue is left on the cell.
$newQ3.Range("B2:G3").NumberFormat = "@"

$newQ3.Range("B2").Value = "004845"
$newQ3.Range("C2").Value = "南华瑞盈混合A"
$newQ3.Range("D2").Value = "2.63"
$newQ3.Range("E2").Value = "93.80"
$newQ3.Range("F2").Value = "2.93"
$newQ3.Range("G2").Value = "0.0771"
$newQ3.Range("H2").Value = 10

$newQ3.Range("A3").Value = 1
$newQ3.Range("B3").Value = "004846"
$newQ3.Range("C3").Value = "南华瑞盈混合C"
$newQ3.Range("D3").Value = "0.08"
$newQ3.Range("E3").Value = "93.80"
$newQ3.Range("F3").Value = "2.93"
$newQ3.Range("G3").Value = "0.0023"
$newQ3.Range("H3").Value = 10

$newQ3.Range("B2:G3").Style = "Normal"

# --- 3. Update the "总计" sheet: insert a new row 2 for 2022-Q3, pushing the
#        existing 2021-Q3 row down to row 3, and renumber the index column.
$total.Rows.Item(2).Insert()
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.08

$total.Range("A3").Value = 1
